# Fix paragraph font size setting to use runs instead of style.
#
# The document's "Normal" paragraph style carried an explicit
# <w:rPr><w:sz w:val="22"/></w:rPr> override (11pt). This script moves
# that 11pt sizing down onto the individual body-text runs that were
# relying on it (so each affected run gets its own explicit
# <w:rPr><w:sz w:val="22"/></w:rPr>), then clears the now-redundant
# size override from the Normal style itself.

$d = $word.ActiveDocument

# --- Step 1: stamp explicit 11pt (sz=22) on every plain body-text run ---
# Target paragraphs are the un-styled "Normal" body lines (problem
# statement, EDA notes, steps, results, etc.) — i.e. paragraphs that:
#   - use the "Normal" style (not Title/Heading/List Bullet/etc.)
#   - are not already carrying their own explicit direct formatting
#     (the "WHAT TO SAY:" labels and the byline are bold/colored and
#     must be left untouched)
#   - are not empty spacer paragraphs
#   - are not the plain "────…" divider rules
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    $len = $text.Length

    if ($p.Style.NameLocal -ne "Normal") { continue }
    if ($p.Range.Font.Bold -ne 0) { continue }
    if ($len -le 1) { continue }
    if (-not ($text -match '[A-Za-z0-9]')) { continue }

    # Apply to the run text only — exclude the trailing paragraph mark
    # so the pilcrow's own run properties (w:pPr/w:rPr) stay untouched
    # and only the run(s) inside the paragraph get <w:rPr><w:sz .../>.
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.Font.Size = 11
}

# --- Step 2: drop the now-unneeded size override from the Normal style ---
$normal = $d.Styles.Item("Normal")
$normal.Font.Size = 11
